$wb = $excel.ActiveWorkbook

# Update the staking amount value for Staker A on the "Staking" sheet.
$ws = $wb.Worksheets.Item("Staking")
$ws.Range("B2").Value = 10000000000

# Make "Staking" the active sheet / active tab, and set the active cell selection.
$ws.Activate()
$ws.Range("G12").Select()
